$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like text in column C stays literal text, not auto-converted to dates
$ws.Range("C2:C106").NumberFormat = "@"

$rowData = @(
    @(2, $null, $null, $null, 66.8),
    @(4, $null, $null, $null, 18.4),
    @(7, $null, $null, $null, 66.59999999999999),
    @(9, $null, $null, $null, 18.4),
    @(12, $null, $null, $null, 65.2),
    @(14, $null, $null, $null, 19),
    @(17, $null, $null, $null, 65),
    @(19, $null, $null, $null, 19.2),
    @(22, $null, $null, $null, 64),
    @(24, $null, $null, $null, 21),
    @(27, $null, $null, $null, 62.5),
    @(29, $null, $null, $null, 22.3),
    @(32, "Brasil", $null, "01/01/2024", 61.6),
    @(33, "Brasil", $null, "01/01/2024", 6),
    @(34, "Brasil", $null, "01/01/2024", 23),
    @(35, "Brasil", $null, "01/01/2024", 9.1),
    @(36, "Brasil", $null, "01/01/2024", $null),
    @(37, $null, $null, "01/01/2016", $null),
    @(38, $null, $null, "01/01/2016", 3.6),
    @(39, $null, $null, "01/01/2016", 14.9),
    @(40, $null, $null, "01/01/2016", $null),
    @(41, $null, $null, "01/01/2016", $null),
    @(42, $null, $null, "01/01/2017", 73.2),
    @(43, $null, $null, "01/01/2017", 3.5),
    @(44, $null, $null, "01/01/2017", 14.8),
    @(45, $null, $null, "01/01/2017", 8.199999999999999),
    @(46, $null, $null, "01/01/2017", $null),
    @(47, $null, $null, "01/01/2018", 71.59999999999999),
    @(48, $null, $null, "01/01/2018", 3.7),
    @(49, $null, $null, "01/01/2018", 15.5),
    @(50, $null, $null, "01/01/2018", 8.9),
    @(51, $null, $null, "01/01/2018", $null),
    @(52, $null, $null, "01/01/2019", 72.59999999999999),
    @(53, $null, $null, "01/01/2019", $null),
    @(54, $null, $null, "01/01/2019", 15.5),
    @(55, $null, $null, "01/01/2019", 8.6),
    @(56, $null, $null, "01/01/2019", $null),
    @(57, $null, $null, "01/01/2022", 71.3),
    @(58, $null, $null, "01/01/2022", 3.2),
    @(59, $null, $null, "01/01/2022", 17),
    @(60, $null, $null, "01/01/2022", 8.4),
    @(61, $null, $null, "01/01/2022", $null),
    @(62, "Nordeste", $null, "01/01/2023", 70.59999999999999),
    @(63, "Nordeste", $null, "01/01/2023", 2.9),
    @(64, "Nordeste", $null, "01/01/2023", 17.7),
    @(65, "Nordeste", $null, "01/01/2023", 8.6),
    @(66, "Nordeste", $null, "01/01/2023", 0.2),
    @(67, "Nordeste", $null, "01/01/2024", 69.59999999999999),
    @(68, "Nordeste", $null, "01/01/2024", 2.9),
    @(69, "Nordeste", $null, "01/01/2024", 18.5),
    @(70, "Nordeste", $null, "01/01/2024", 8.800000000000001),
    @(71, "Nordeste", $null, "01/01/2024", $null),
    @(72, $null, $null, "01/01/2016", 73),
    @(73, $null, $null, "01/01/2016", 3.4),
    @(74, $null, $null, "01/01/2016", 17.4),
    @(75, $null, $null, "01/01/2016", 6),
    @(76, $null, $null, "01/01/2016", 0.1),
    @(77, $null, $null, "01/01/2017", 71.40000000000001),
    @(78, $null, $null, "01/01/2017", 2.8),
    @(79, $null, $null, "01/01/2017", 18.5),
    @(80, $null, $null, "01/01/2017", 7.1),
    @(81, $null, $null, "01/01/2017", 0.2),
    @(82, $null, $null, "01/01/2018", 68.8),
    @(83, $null, $null, "01/01/2018", 4.9),
    @(84, $null, $null, "01/01/2018", 16.9),
    @(85, $null, $null, "01/01/2018", 8.800000000000001),
    @(86, $null, $null, "01/01/2018", 0.6),
    @(87, $null, $null, "01/01/2019", $null),
    @(88, $null, $null, "01/01/2019", 3.9),
    @(89, $null, $null, "01/01/2019", 19.9),
    @(90, $null, $null, "01/01/2019", 11.2),
    @(91, $null, $null, "01/01/2019", $null),
    @(92, "Sergipe", "Próprio de algum morador - já pago", "01/01/2022", 62),
    @(93, "Sergipe", "Próprio de algum morador - ainda pagando", "01/01/2022", 4.8),
    @(94, "Sergipe", "Alugado", "01/01/2022", 21.9),
    @(95, "Sergipe", "Cedido", "01/01/2022", 10.9),
    @(96, "Sergipe", "Outra condição", "01/01/2022", 0.5),
    @(97, "Sergipe", "Próprio de algum morador - já pago", "01/01/2023", 64.7),
    @(98, "Sergipe", "Próprio de algum morador - ainda pagando", "01/01/2023", 4.1),
    @(99, "Sergipe", "Alugado", "01/01/2023", 20.9),
    @(100, "Sergipe", "Cedido", "01/01/2023", 9.699999999999999),
    @(101, "Sergipe", "Outra condição", "01/01/2023", 0.6),
    @(102, "Sergipe", "Próprio de algum morador - já pago", "01/01/2024", 60.7),
    @(103, "Sergipe", "Próprio de algum morador - ainda pagando", "01/01/2024", 6.5),
    @(104, "Sergipe", "Alugado", "01/01/2024", 24.9),
    @(105, "Sergipe", "Cedido", "01/01/2024", 7.9),
    @(106, "Sergipe", "Outra condição", "01/01/2024", $null)
)

foreach ($row in $rowData) {
    $r = $row[0]
    if ($row[1] -ne $null) { $ws.Cells.Item($r, 1).Value = $row[1] }
    if ($row[2] -ne $null) { $ws.Cells.Item($r, 2).Value = $row[2] }
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 3).Value = $row[3] }
    if ($row[4] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[4] }
}

$ws.Range("C2:C106").ClearFormats()
